$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 382-383, pushing the existing data (old rows 382-410)
# down to rows 384-412. This matches the way the source data was extended
# with a newer weekly price report while keeping the historical rows intact.
$ws.Rows("382:383").Insert()

# Row 382: new "Primera" entry for date serial 44706 (2022-05-25)
$ws.Cells.Item(382, 1).Value = 8
$ws.Cells.Item(382, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(382, 3).Value = "Coquimbo"
$ws.Cells.Item(382, 4).Value = 44706
$ws.Cells.Item(382, 5).Value = 4
$ws.Cells.Item(382, 6).Value = 100112043
$ws.Cells.Item(382, 7).Value = "Pepino dulce"
$ws.Cells.Item(382, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 560
$ws.Cells.Item(382, 11).Value = 13000
$ws.Cells.Item(382, 12).Value = 14000
$ws.Cells.Item(382, 13).Value = 13500
$ws.Cells.Item(382, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(382, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(382, 16).Value = 750
$ws.Cells.Item(382, 17).Value = 18
$ws.Cells.Item(382, 18).Value = "Hortaliza"

# Row 383: new "Segunda" entry for date serial 44706 (2022-05-25)
$ws.Cells.Item(383, 1).Value = 8
$ws.Cells.Item(383, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(383, 3).Value = "Coquimbo"
$ws.Cells.Item(383, 4).Value = 44706
$ws.Cells.Item(383, 5).Value = 4
$ws.Cells.Item(383, 6).Value = 100112043
$ws.Cells.Item(383, 7).Value = "Pepino dulce"
$ws.Cells.Item(383, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(383, 9).Value = "Segunda"
$ws.Cells.Item(383, 10).Value = 400
$ws.Cells.Item(383, 11).Value = 10000
$ws.Cells.Item(383, 12).Value = 11000
$ws.Cells.Item(383, 13).Value = 10500
$ws.Cells.Item(383, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(383, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(383, 16).Value = 583
$ws.Cells.Item(383, 17).Value = 18
$ws.Cells.Item(383, 18).Value = "Hortaliza"
